$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON): clear B2 entirely, update C2:E2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.2558677754913283
$ws.Range("D2").Value = 0.15274354883031485
$ws.Range("E2").Value = 3.1094321911993568

# Row 3 (STR): update B3:E3
$ws.Range("B3").Value = 0.2347056316217849
$ws.Range("C3").Value = 4.678302885769237
$ws.Range("D3").Value = 0.28569769304332177
$ws.Range("E3").Value = 2.8182186993434706

# Update the selection shown when the sheet is opened
$ws.Range("B1:E3").Select()
